$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Burndown "real" series (row 4): fix R4 and extend with S4:U4 ---
$ws.Range("R4").Value = 42
$ws.Range("S4").Value = 42
$ws.Range("T4").Value = 38
$ws.Range("U4").Value = 36

# --- "finished" counts (row 6): fix Q6 and extend with R6:T6 ---
$ws.Range("Q6").Value = 22
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 2

# --- Actual-hours breakdown table (rows 11-19): add week 2 (D) / week 3 (E)
#     columns and per-row "sprint 1" totals (F), then refresh the column totals
#     in row 20 ---
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 5
$ws.Range("F11").Formula = "=SUM(C11,D11,E11)"

$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 10

$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 0

$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0

$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 0

$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 0

$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 1

$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 2

$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 8

# Fill F12:F19 together so the engine stores it as one shared formula group
# (matches row 12's master formula + si="0" follower cells on rows 13-19).
$ws.Range("F12:F19").Formula = "=SUM(C12,D12,E12)"

$ws.Range("D20").Formula = "=SUM(D11:D19)"
$ws.Range("E20").Formula = "=SUM(E11:E19)"
$ws.Range("F20").Formula = "=SUM(F11:F19)"

# --- Chart: "real" series now only runs through column U (matches the
#     widened row 4 range above) ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$realSeries = $chart.SeriesCollection().Item(2)
$realSeries.Formula = "=SERIES(Sheet1!`$A`$4,Sheet1!`$B`$2:`$W`$2,Sheet1!`$B`$4:`$U`$4,2)"

# --- Page setup: explicit portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection cursor moved to F26 ahead of the end-of-sprint meeting ---
$ws.Range("F26").Select()
